$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, reusing the same header style as
# the rest of row 1 (bold, centered, bordered) by copying formats from E1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Per-row timestamps for the new column (plain, unstyled cells like the
# rest of the data rows).
$times = @(
    "2021-10-05 10:52:48.204189",
    "2021-10-05 10:52:48.204202",
    "2021-10-05 10:52:48.204206",
    "2021-10-05 10:52:48.204209",
    "2021-10-05 10:52:48.204212",
    "2021-10-05 10:52:48.204216",
    "2021-10-05 10:52:48.204219"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
